$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data in columns D (Price) and E (Volume(1h)) is stored as
# plain text (e.g. "332.82", "0.95%"), not as numbers/percentages.
# Force the Text number format on the cells we are about to rewrite so
# Excel keeps storing them as text instead of auto-converting them to
# numeric values.
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D9:D15").NumberFormat = "@"
$ws.Range("D17:D25").NumberFormat = "@"
$ws.Range("D39:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"
$ws.Range("E39:E51").NumberFormat = "@"

# Apply the updated cell values from the diff.
$ws.Range("D2").Value = "332.82"
$ws.Range("E2").Value = "0.95%"
$ws.Range("D3").Value = "41.34"
$ws.Range("E3").Value = "2.26%"
$ws.Range("D4").Value = "5.727"
$ws.Range("E4").Value = "-3.58%"
$ws.Range("D5").Value = "0.08111"
$ws.Range("E5").Value = "-0.44%"
$ws.Range("D6").Value = "2.083"
$ws.Range("E6").Value = "6.57%"
$ws.Range("D7").Value = "8.745"
$ws.Range("E7").Value = "-0.03%"
$ws.Range("E8").Value = "-0.84%"
$ws.Range("D9").Value = "2.960"
$ws.Range("E9").Value = "-1.30%"
$ws.Range("D10").Value = "0.9253"
$ws.Range("E10").Value = "-1.97%"
$ws.Range("D11").Value = "0.1280"
$ws.Range("E11").Value = "-2.18%"
$ws.Range("D12").Value = "0.1967"
$ws.Range("E12").Value = "-1.58%"
$ws.Range("D13").Value = "8.807"
$ws.Range("E13").Value = "14.10%"
$ws.Range("D14").Value = "0.09294"
$ws.Range("E14").Value = "0.53%"
$ws.Range("D15").Value = "0.03717"
$ws.Range("E15").Value = "8.32%"
$ws.Range("E16").Value = "9.22%"
$ws.Range("D17").Value = "0.001314"
$ws.Range("E17").Value = "-1.32%"
$ws.Range("D18").Value = "0.006130"
$ws.Range("E18").Value = "1.85%"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "0.004436"
$ws.Range("E19").Value = "1.79%"
$ws.Range("B20").Value = "LEO"
$ws.Range("C20").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D20").Value = "3.380"
$ws.Range("E20").Value = "0.20%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "0.3529"
$ws.Range("E21").Value = "0.93%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "0.1417"
$ws.Range("E22").Value = "-1.60%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2610"
$ws.Range("E23").Value = "6.53%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "0.04440"
$ws.Range("E24").Value = "0.45%"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "0.001259"
$ws.Range("E25").Value = "0.47%"
$ws.Range("E26").Value = "4.40%"
$ws.Range("D39").Value = "0.02840"
$ws.Range("E39").Value = "13.41%"
$ws.Range("D40").Value = "0.05512"
$ws.Range("E40").Value = "3.84%"
$ws.Range("D41").Value = "0.007711"
$ws.Range("E41").Value = "1.46%"
$ws.Range("D42").Value = "0.009914"
$ws.Range("E42").Value = "11.13%"
$ws.Range("D43").Value = "0.1423"
$ws.Range("E43").Value = "-0.69%"
$ws.Range("D44").Value = "0.002090"
$ws.Range("E44").Value = "1.43%"
$ws.Range("D45").Value = "0.01100"
$ws.Range("E45").Value = "4.89%"
$ws.Range("D46").Value = "0.00006789"
$ws.Range("E46").Value = "-1.10%"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48").Value = "0.002992"
$ws.Range("E48").Value = "3.42%"
$ws.Range("D49").Value = "0.002279"
$ws.Range("E49").Value = "26.56%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.05%"
